# Insert a new data row at row 289 (shifting the existing rows 289-326
# down to 290-327) and populate it with the new weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(289).Insert()

$ws.Range("A289").Value2 = 6
$ws.Range("B289").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C289").Value2 = "Metropolitana"
$ws.Range("D289").Value2 = 45131
$ws.Range("E289").Value2 = 13
$ws.Range("F289").Value2 = 100112022
$ws.Range("G289").Value2 = "Arveja Verde"
$ws.Range("H289").Value2 = "Perfection"
$ws.Range("I289").Value2 = "Primera"
$ws.Range("J289").Value2 = 200
$ws.Range("K289").Value2 = 37000
$ws.Range("L289").Value2 = 38000
$ws.Range("M289").Value2 = 37500
$ws.Range("N289").Value2 = "`$/malla 25 kilos"
$ws.Range("O289").Value2 = "Provincia de Huasco"
$ws.Range("P289").Value2 = 1500
$ws.Range("Q289").Value2 = 25
$ws.Range("R289").Value2 = "Hortaliza"

# Make sure the D289 cell keeps the date-formatted style (s="2") like the
# rest of the column, and that the date serial renders as an actual date.
$ws.Range("D289").NumberFormat = $ws.Range("D290").NumberFormat
